$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Globo"
$ws.Range("B9").Value = "RJ TV 1"
$ws.Range("C9").Value = "Social"
$ws.Range("D9").Value = "2025-04-01T12:36"
$ws.Range("E9").Value = "Positivo"
$ws.Range("F9").Value = "Oportunidades de trabalho. Em Campos, 366 vagas, entre elas para taifeiro e nutricionista offshore, Garçom e auxiliar de serviços gerais.  "
